# [ADDITIONAL SCRAPING] added scraping code for extra bowling attributes and
# excel sheets.
#
# Adds a new "ODI Bowling Extra" worksheet (mirroring the existing
# "ODI Batting Extra" sheet) after the last sheet in the workbook, and fills
# it in with the scraped MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL stats per
# MATCH_CODE.

$wb = $excel.ActiveWorkbook

# Add the new worksheet as the last sheet in the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ODI Bowling Extra"

# Header row: same 3 headers as other "extra" scraping sheets.
$newSheet.Cells.Item(1, 1).Value = "MATCH_CODE"
$newSheet.Cells.Item(1, 2).Value = "MAIDEN_OVERS"
$newSheet.Cells.Item(1, 3).Value = "PERCENT_WICKETS_OF_ALL"

# Re-use the bold / bordered / centered header style already used on the
# other scraping sheets (e.g. "ODI Batting Extra" row 1) instead of building
# a brand-new style.
$existingHeader = $wb.Worksheets.Item("ODI Batting Extra").Range("A1")
$existingHeader.Copy()
$newSheet.Range("A1:C1").PasteSpecial(-4122)

# MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL values, scraped per
# match. Every value (even "0") is stored as text, matching how the other
# scraping sheets store their columns; blank entries are left as empty text
# cells rather than dropped.
$data = @(
    @("4300", "0", ""),
    @("4303", "", ""),
    @("4308", "", ""),
    @("4314", "1", "30.00%"),
    @("4321", "0", ""),
    @("4326", "", ""),
    @("4331", "0", ""),
    @("4336", "0", "10.00%"),
    @("4342", "0", ""),
    @("4346", "", ""),
    @("4354", "0", ""),
    @("4355", "", ""),
    @("4454", "1", "30.00%"),
    @("4456", "", ""),
    @("4457", "0", "10.00%"),
    @("4472", "0", ""),
    @("4476", "0", ""),
    @("4609", "0", ""),
    @("4618", "0", ""),
    @("4619", "", "")
)

$r = 2
foreach ($row in $data) {
    for ($c = 1; $c -le 3; $c++) {
        $cell = $newSheet.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$c - 1]
        $cell.ClearFormats()
    }
    $r = $r + 1
}
